$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 251562
$ws.Range("C2").Value = 532.9449956670721
$ws.Range("D2").Value = 113.0475835137304
$ws.Range("E2").Value = 392
$ws.Range("F2").Value = 453
$ws.Range("G2").Value = 498
$ws.Range("H2").Value = 577
$ws.Range("I2").Value = 1593

$ws.Range("B3").Value = 251562
$ws.Range("C3").Value = 44.98589846638205
$ws.Range("D3").Value = 4.890338910421502
$ws.Range("E3").Value = 30.48
$ws.Range("F3").Value = 41.31
$ws.Range("G3").Value = 44.73
$ws.Range("H3").Value = 48.4
$ws.Range("I3").Value = 60.19

$ws.Range("B4").Value = 251562
$ws.Range("C4").Value = 1.431533816713176
$ws.Range("D4").Value = 3.591974203725274
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.52
$ws.Range("G4").Value = 1.06
$ws.Range("H4").Value = 1.92
$ws.Range("I4").Value = 637.71

$ws.Range("B5").Value = 251562
$ws.Range("C5").Value = 319.7001681096509
$ws.Range("D5").Value = 9.851875056702589
$ws.Range("E5").Value = 286.91
$ws.Range("F5").Value = 313.79
$ws.Range("G5").Value = 321.37
$ws.Range("H5").Value = 325.75
$ws.Range("I5").Value = 342

$ws.Range("B6").Value = 251562
$ws.Range("C6").Value = 22.58832927866689
$ws.Range("D6").Value = 1.9432519255438
$ws.Range("E6").Value = 15.83
$ws.Range("F6").Value = 21.5
$ws.Range("G6").Value = 22.2
$ws.Range("H6").Value = 23.49
$ws.Range("I6").Value = 33.93

$ws.Range("B7").Value = 251562
$ws.Range("C7").Value = -76.56134471820069
$ws.Range("D7").Value = 23.72077317696776
$ws.Range("E7").Value = -123
$ws.Range("F7").Value = -94
$ws.Range("G7").Value = -74
$ws.Range("H7").Value = -54
$ws.Range("I7").Value = -32

$ws.Range("B8").Value = 251059
$ws.Range("C8").Value = 7.671104800066916
$ws.Range("D8").Value = 6.523565966186341
$ws.Range("E8").Value = -24.5
$ws.Range("F8").Value = 7.8
$ws.Range("G8").Value = 9.5
$ws.Range("H8").Value = 11.2
$ws.Range("I8").Value = 19

$ws.Range("B9").Value = 251562
$ws.Range("C9").Value = 9.324361390035062
$ws.Range("D9").Value = 1.688266393300287
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = 11
$ws.Range("I9").Value = 12

$ws.Range("B10").Value = 251562
$ws.Range("C10").Value = 867.8306803094266
$ws.Range("D10").Value = 0.4611932605059126
$ws.Range("E10").Value = 867.1
$ws.Range("F10").Value = 867.5
$ws.Range("G10").Value = 867.9
$ws.Range("H10").Value = 868.3
$ws.Range("I10").Value = 868.5

$ws.Range("B11").Value = 251550
$ws.Range("C11").Value = 17213.86838799443
$ws.Range("D11").Value = 11936.9850653751
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 6046.25
$ws.Range("G11").Value = 16487.5
$ws.Range("H11").Value = 27836
$ws.Range("I11").Value = 40744

$ws.Range("B12").Value = 251562
$ws.Range("C12").Value = 19240.49235973637
$ws.Range("D12").Value = 13294.66376833996
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 6817
$ws.Range("G12").Value = 18470
$ws.Range("H12").Value = 31111
$ws.Range("I12").Value = 44461

$ws.Range("B13").Value = 251562
$ws.Range("C13").Value = 0.5573166730110272
$ws.Range("D13").Value = 0.5906980832801664
$ws.Range("E13").Value = 0.071936
$ws.Range("F13").Value = 0.133632
$ws.Range("G13").Value = 0.246784
$ws.Range("H13").Value = 0.987136
$ws.Range("I13").Value = 1.974272

$ws.Range("B14").Value = 251562
$ws.Range("C14").Value = 23.89372798753389
$ws.Range("D14").Value = 13.4170136097565
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = 39
$ws.Range("I14").Value = 43

$ws.Range("B15").Value = 251562
$ws.Range("C15").Value = 0.6732614623830309
$ws.Range("D15").Value = 0.7487418961972424
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 2

$ws.Range("B16").Value = 251562
$ws.Range("C16").Value = 1.830041898219922
$ws.Range("D16").Value = 1.670023765759874
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 5

$ws.Range("B17").Value = 251562
$ws.Range("C17").Value = 93.96134471820075
$ws.Range("D17").Value = 23.72077317695606
$ws.Range("E17").Value = 49.4
$ws.Range("F17").Value = 71.4
$ws.Range("G17").Value = 91.4
$ws.Range("H17").Value = 111.4
$ws.Range("I17").Value = 140.4

$ws.Range("B18").Value = 251059
$ws.Range("C18").Value = -85.67615448186591
$ws.Range("D18").Value = 21.4250801593023
$ws.Range("E18").Value = -125.9574620641016
$ws.Range("F18").Value = -103.3377954106368
$ws.Range("G18").Value = -85.26572375596102
$ws.Range("H18").Value = -65.14699179957641
$ws.Range("I18").Value = -33.14609373817283

$ws.Range("B19").Value = 251059
$ws.Range("C19").Value = -78.00504968179897
$ws.Range("D19").Value = 25.73412961331832
$ws.Range("E19").Value = -145.0217119216414
$ws.Range("F19").Value = -94.26572375596102
$ws.Range("G19").Value = -74.39612087980606
$ws.Range("H19").Value = -54.34699179957641
$ws.Range("I19").Value = -33.49305820175223

